$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to remain plain text (preserving
# the exact literal string, e.g. "0.5067" or "1.010") and without leaving any
# permanent number-format / style change behind on the cell.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '28.465.13'
Set-TextValue 'E2' '  -0.21%  '
Set-TextValue 'D3' '1.866.49'
Set-TextValue 'E3' '  -0.53%  '
Set-TextValue 'E4' '  -1.68%  '
Set-TextValue 'E5' '  -0.93%  '
Set-TextValue 'E6' '  -1.40%  '
Set-TextValue 'D7' '0.5067'
Set-TextValue 'E7' '  -1.69%  '
Set-TextValue 'D8' '0.3907'
Set-TextValue 'E8' '  -1.94%  '
Set-TextValue 'D9' '0.08321'
Set-TextValue 'E9' '  -0.92%  '
Set-TextValue 'D10' '42.58'
Set-TextValue 'E10' '  +0.95%  '
Set-TextValue 'D11' '1.102'
Set-TextValue 'E11' '  -1.19%  '
Set-TextValue 'D12' '6.190'
Set-TextValue 'E12' '  -1.34%  '
Set-TextValue 'D13' '1.868.70'
Set-TextValue 'E13' '  +2.27%  '
Set-TextValue 'D14' '20.29'
Set-TextValue 'E14' '  -1.36%  '
Set-TextValue 'E15' '  -0.26%  '
Set-TextValue 'D16' '1.010'
Set-TextValue 'E16' '  -1.66%  '
Set-TextValue 'D17' '0.00001098'
Set-TextValue 'E17' '  -1.34%  '
Set-TextValue 'D18' '91.19'
Set-TextValue 'E18' '  -0.13%  '
Set-TextValue 'D19' '0.06728'
Set-TextValue 'E19' '  -0.87%  '
Set-TextValue 'D20' '17.62'
Set-TextValue 'E20' '  -0.92%  '
Set-TextValue 'D21' '1.007'
Set-TextValue 'E21' '  -1.44%  '
Set-TextValue 'D22' '5.904'
Set-TextValue 'E22' '  -1.31%  '
Set-TextValue 'D23' '28.516.13'
Set-TextValue 'E23' '  -0.04%  '
Set-TextValue 'D24' '11.06'
Set-TextValue 'E24' '  -1.24%  '
Set-TextValue 'E25' '  -3.74%  '
Set-TextValue 'D26' '2.078.10'
Set-TextValue 'E26' '  +1.99%  '
Set-TextValue 'D27' '157.69'
Set-TextValue 'E27' '  -2.90%  '
Set-TextValue 'D28' '20.53'
Set-TextValue 'E28' '  -1.25%  '
Set-TextValue 'D29' '2.411'
Set-TextValue 'E29' '  +1.59%  '
Set-TextValue 'D30' '125.63'
Set-TextValue 'E30' '  -1.78%  '
Set-TextValue 'E31' '  -1.59%  '
Set-TextValue 'D32' '1.034'
Set-TextValue 'E32' '  -0.67%  '
Set-TextValue 'D33' '5.741'
Set-TextValue 'E33' '  -1.54%  '
Set-TextValue 'D34' '3.619'
Set-TextValue 'E34' '  -0.78%  '
Set-TextValue 'B35' 'VeChain'
Set-TextValue 'C35' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D35' '0.02441'
Set-TextValue 'E35' '  -0.06%  '
Set-TextValue 'B36' 'Hedera'
Set-TextValue 'C36' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D36' '0.06584'
Set-TextValue 'E36' '  +1.08%  '
Set-TextValue 'E37' '  +0.19%  '
Set-TextValue 'D38' '0.2154'
Set-TextValue 'E38' '  -1.70%  '
Set-TextValue 'D39' '5.014'
Set-TextValue 'E39' '  -0.66%  '
Set-TextValue 'D40' '1.177'
Set-TextValue 'E40' '  -1.27%  '
Set-TextValue 'D41' '1.233'
Set-TextValue 'E41' '  -4.04%  '
Set-TextValue 'E42' '  -1.57%  '
Set-TextValue 'E43' '  -1.84%  '
Set-TextValue 'E44' '  -1.17%  '
Set-TextValue 'D45' '0.5988'
Set-TextValue 'E45' '  -0.98%  '
Set-TextValue 'D46' '13.05'
Set-TextValue 'E46' '  -0.71%  '
Set-TextValue 'D47' '3.683'
Set-TextValue 'E47' '  -1.42%  '
Set-TextValue 'D48' '1.993'
Set-TextValue 'E48' '  -0.32%  '
Set-TextValue 'D49' '122.12'
Set-TextValue 'E49' '  -0.19%  '
Set-TextValue 'D50' '1.207'
Set-TextValue 'E50' '  -0.37%  '
Set-TextValue 'D51' '1.132'
Set-TextValue 'E51' '  -7.68%  '
